# Generate Report for Handback
# Updates generated-timestamp and priority values across the Overview,
# zh-cn and de-de sheets of the handback-status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 20:14:32"
$wsOverview.Range("G5").Value = "2016-08-25 20:14:32"

# --- zh-cn sheet ---
# Column E = "Priority", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-25 20:14:27"
$wsZhCn.Range("H5").Value = "2016-08-25 20:14:27"
$wsZhCn.Range("K2").Value = "2016-08-25 20:14:42"
$wsZhCn.Range("K5").Value = "2016-08-25 20:14:42"

# --- de-de sheet ---
# Column E = "Priority", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-25 20:14:32"
$wsDeDe.Range("H5").Value = "2016-08-25 20:14:32"
$wsDeDe.Range("K2").Value = "2016-08-25 20:14:49"
$wsDeDe.Range("K5").Value = "2016-08-25 20:14:49"
